$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$result = $wb.Worksheets.Item("result")

# Insert two new rows at row 2 to make room for the new entries,
# pushing the existing Allocation/Expend/Remaining block down.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# New transaction rows
$ws.Cells.Item(2, 1).Value = "خوراکی"
$ws.Cells.Item(2, 2).Value = 30000
$ws.Cells.Item(2, 3).Value = "'1400/01/01"
$ws.Cells.Item(2, 3).Style = "Normal"

$ws.Cells.Item(3, 1).Value = "خوراکی"
$ws.Cells.Item(3, 2).Value = 10000
$ws.Cells.Item(3, 3).Value = "'1400/01/10"
$ws.Cells.Item(3, 3).Style = "Normal"

# Update Expend (row 6) and Remaining (row 7) after the insert
$ws.Cells.Item(6, 2).Value = 40000
$ws.Cells.Item(7, 2).Value = 160000

# Update the result sheet summary (row 2 corresponds to تفریح)
$result.Cells.Item(2, 3).Value = 40000
$result.Cells.Item(2, 4).Value = 160000

$result.Cells.Item(12, 3).Value = 2040000
$result.Cells.Item(12, 4).Value = 460000
